$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating the existing
#    "2022-Q2" sheet (tab 2). Copying - rather than Worksheets.Add() - means
#    the new tab automatically inherits the matching sheetPr/pageMargins/
#    column-A & header styling used by its sibling quarter sheets, and it is
#    placed directly after "总计" / before "2022-Q2", exactly like the diff.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2template = $wb.Worksheets.Item(2)
$q2template.Copy($null, $totalSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template has 9 data rows (2022-Q2 held 8 funds); 2022-Q3 only has 7
# funds (8 rows incl. header), so drop the spare last row.
$q3.Rows.Item(9).Delete()

# "2022-Q2" (the pristine, untouched template) is now pushed to tab 3 - keep
# a handle on it so we can borrow its cell formatting for the text columns
# below without disturbing the values we are about to overwrite.
$q2 = $wb.Worksheets.Item(3)

# Header row text (values only - formatting/style already correct from the
# sheet copy).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

function Set-TextCell($sheet, $addr, $value, $formatSourceSheet) {
    # Force the cell to Text so numeric-looking strings (fund codes like
    # "013385", decimals like "0.4988") keep their literal representation
    # instead of being coerced into a number, then restore the original
    # (non-"@") number format from the pristine template sheet so the style
    # index matches what a hand-authored file would use.
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $value
    $formatSourceSheet.Range($addr).Copy()
    $sheet.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}

# Row 2
$q3.Range("A2").Value = 0
Set-TextCell $q3 "B2" "013385" $q2
Set-TextCell $q3 "C2" "信澳优势价值混合A" $q2
Set-TextCell $q3 "D2" "12.44" $q2
Set-TextCell $q3 "E2" "84.28" $q2
Set-TextCell $q3 "F2" "4.01" $q2
Set-TextCell $q3 "G2" "0.4988" $q2
$q3.Range("H2").Value = 6

# Row 3
$q3.Range("A3").Value = 1
Set-TextCell $q3 "B3" "004685" $q2
Set-TextCell $q3 "C3" "金元顺安元启灵活配置混合" $q2
Set-TextCell $q3 "D3" "15.28" $q2
Set-TextCell $q3 "E3" "77.14" $q2
Set-TextCell $q3 "F3" "0.85" $q2
Set-TextCell $q3 "G3" "0.1299" $q2
$q3.Range("H3").Value = 10

# Row 4
$q3.Range("A4").Value = 2
Set-TextCell $q3 "B4" "013393" $q2
Set-TextCell $q3 "C4" "信澳价值精选混合A" $q2
Set-TextCell $q3 "D4" "3.34" $q2
Set-TextCell $q3 "E4" "79.98" $q2
Set-TextCell $q3 "F4" "3.76" $q2
Set-TextCell $q3 "G4" "0.1256" $q2
$q3.Range("H4").Value = 4

# Row 5
$q3.Range("A5").Value = 3
Set-TextCell $q3 "B5" "013386" $q2
Set-TextCell $q3 "C5" "信澳优势价值混合C" $q2
Set-TextCell $q3 "D5" "1.26" $q2
Set-TextCell $q3 "E5" "84.28" $q2
Set-TextCell $q3 "F5" "4.01" $q2
Set-TextCell $q3 "G5" "0.0505" $q2
$q3.Range("H5").Value = 6

# Row 6
$q3.Range("A6").Value = 4
Set-TextCell $q3 "B6" "013394" $q2
Set-TextCell $q3 "C6" "信澳价值精选混合C" $q2
Set-TextCell $q3 "D6" "0.38" $q2
Set-TextCell $q3 "E6" "79.98" $q2
Set-TextCell $q3 "F6" "3.76" $q2
Set-TextCell $q3 "G6" "0.0143" $q2
$q3.Range("H6").Value = 4

# Row 7
$q3.Range("A7").Value = 5
Set-TextCell $q3 "B7" "164811" $q2
Set-TextCell $q3 "C7" "工银瑞信中证京津冀协同发展主题指数（LOF）A" $q2
Set-TextCell $q3 "D7" "0.12" $q2
Set-TextCell $q3 "E7" "93.09" $q2
Set-TextCell $q3 "F7" "2.95" $q2
Set-TextCell $q3 "G7" "0.0035" $q2
$q3.Range("H7").Value = 9

# Row 8
$q3.Range("A8").Value = 6
Set-TextCell $q3 "B8" "164825" $q2
Set-TextCell $q3 "C8" "工银瑞信中证京津冀协同发展主题指数（LOF）C" $q2
Set-TextCell $q3 "D8" "0.03" $q2
Set-TextCell $q3 "E8" "93.09" $q2
Set-TextCell $q3 "F8" "2.95" $q2
Set-TextCell $q3 "G8" "0.0009" $q2
$q3.Range("H8").Value = 9

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add a 2022-Q3 row and push the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 4
$totalSheet.Range("D5").Value = 0.05

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 6
$totalSheet.Range("D4").Value = 0.33

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 8
$totalSheet.Range("D3").Value = 0.92

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.82

# A5 is a brand-new cell (the sheet used to stop at row 4) so it needs the
# same bold/border/center styling the other index cells in column A use.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$totalSheet.Activate()
